# khl/Injuries_Master_Clubs.xlsx refresh (2025-11-26 03:0x UTC scrape)
# - 5 players recovered (left "snapshot" of currently-injured players)
# - those same 5 players are appended to the "returned" log for 2025-11-26
# - every remaining "snapshot" row gets a fresh scraped_at timestamp

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) snapshot: drop the players who returned from injury.
#    Delete bottom-to-top so earlier row numbers stay valid.
# ---------------------------------------------------------------
$snapshot = $wb.Worksheets.Item("snapshot")

$returnedRows = 37, 33, 31, 30, 26
foreach ($r in $returnedRows) {
    $snapshot.Rows.Item($r).Delete()
}

# ---------------------------------------------------------------
# 2) snapshot: refresh scraped_at (column K) for every remaining
#    data row (2..34) to the new scrape run's timestamps.
# ---------------------------------------------------------------
$newTimestamps = @(
    "2025-11-26T03:02:55.353047+00:00",
    "2025-11-26T03:02:57.632335+00:00",
    "2025-11-26T03:02:57.632369+00:00",
    "2025-11-26T03:03:00.428256+00:00",
    "2025-11-26T03:03:02.776613+00:00",
    "2025-11-26T03:03:05.108922+00:00",
    "2025-11-26T03:03:05.108951+00:00",
    "2025-11-26T03:03:05.108969+00:00",
    "2025-11-26T03:03:07.338385+00:00",
    "2025-11-26T03:03:10.063262+00:00",
    "2025-11-26T03:03:12.738729+00:00",
    "2025-11-26T03:03:15.456918+00:00",
    "2025-11-26T03:03:18.209923+00:00",
    "2025-11-26T03:03:23.199849+00:00",
    "2025-11-26T03:03:23.199875+00:00",
    "2025-11-26T03:03:23.199897+00:00",
    "2025-11-26T03:03:25.977597+00:00",
    "2025-11-26T03:03:25.977627+00:00",
    "2025-11-26T03:03:25.977645+00:00",
    "2025-11-26T03:03:25.977661+00:00",
    "2025-11-26T03:03:28.293905+00:00",
    "2025-11-26T03:03:28.293937+00:00",
    "2025-11-26T03:03:31.027453+00:00",
    "2025-11-26T03:03:31.027484+00:00",
    "2025-11-26T03:03:31.027502+00:00",
    "2025-11-26T03:03:31.027522+00:00",
    "2025-11-26T03:03:33.803880+00:00",
    "2025-11-26T03:03:36.088681+00:00",
    "2025-11-26T03:03:38.931494+00:00",
    "2025-11-26T03:03:38.931526+00:00",
    "2025-11-26T03:03:44.044564+00:00",
    "2025-11-26T03:03:46.754572+00:00",
    "2025-11-26T03:03:46.754603+00:00"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $snapshot.Cells.Item($row, 11).Value = $newTimestamps[$i]
}

# ---------------------------------------------------------------
# 3) returned: replace the previous day's log with today's batch
#    of players who came back from injury.
# ---------------------------------------------------------------
$returned = $wb.Worksheets.Item("returned")

$returned.Rows.Item(3).Delete()
$returned.Rows.Item(2).Delete()

$changedAt = "2025-11-26T11:03:47.264605+08:00"
$changedDay = "2025-11-26"

$returnedPlayers = @(
    @("СОЧ", "ХК Сочи",       "Самсонов Илья",   "1369_СОЧ_самсоновилья"),
    @("СЮЛ", "Салават Юлаев", "Берлёв Антон",    "1369_СЮЛ_берлевантон"),
    @("СЮЛ", "Салават Юлаев", "Зоркин Никита",   "1369_СЮЛ_зоркинникита"),
    @("СЮЛ", "Салават Юлаев", "Ян Денис",        "1369_СЮЛ_янденис"),
    @("ЦСК", "ЦСКА",          "Уильямс Колби",   "1369_ЦСК_уильямсколби")
)

for ($i = 0; $i -lt $returnedPlayers.Length; $i++) {
    $row = $i + 2
    $p = $returnedPlayers[$i]
    $returned.Cells.Item($row, 1).Value = $p[0]
    $returned.Cells.Item($row, 2).Value = $p[1]
    $returned.Cells.Item($row, 3).Value = $p[2]
    $returned.Cells.Item($row, 4).Value = $p[3]
    $returned.Cells.Item($row, 5).Value = "RETURN"
    $returned.Cells.Item($row, 6).Value = $changedAt
    # Force plain text so the ISO date string ("2025-11-26") isn't
    # auto-converted into a date serial number, then drop the
    # temporary text format so the cell keeps the workbook's default
    # (unstyled) look, matching the rest of the sheet.
    $returned.Cells.Item($row, 7).NumberFormat = "@"
    $returned.Cells.Item($row, 7).Value = $changedDay
    $returned.Cells.Item($row, 7).ClearFormats()
}
